# ------------------------------------------------------------------
# Reorders rows 2-43 on the "template" sheet to match the updated
# scenario-table row order (see commit "update order of rows in
# template"). Row 1 (header) and rows 44-46 (glucose trace) keep
# their position; only the settings block (2-15) and the scenario
# block (16-43) are permuted. We stage each source row far below the
# used range (offset +1000) to avoid clobbering rows we still need to
# read, then cut everything back into its new home in one pass.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("template")

# --- Phase 1: stage current rows 2-43 down at rows 1002-1043 ---
$ws.Range("A2:B2").Copy($ws.Range("A1002:B1002"))
$ws.Range("A3:B3").Copy($ws.Range("A1003:B1003"))
$ws.Range("A4:B4").Copy($ws.Range("A1004:B1004"))
$ws.Range("A5:B5").Copy($ws.Range("A1005:B1005"))
$ws.Range("A6:B6").Copy($ws.Range("A1006:B1006"))
$ws.Range("A7:B7").Copy($ws.Range("A1007:B1007"))
$ws.Range("A8:B8").Copy($ws.Range("A1008:B1008"))
$ws.Range("A9:B9").Copy($ws.Range("A1009:B1009"))
$ws.Range("A10:B10").Copy($ws.Range("A1010:B1010"))
$ws.Range("A11:B11").Copy($ws.Range("A1011:B1011"))
$ws.Range("A12:B12").Copy($ws.Range("A1012:B1012"))
$ws.Range("A13:B13").Copy($ws.Range("A1013:B1013"))
$ws.Range("A14:B14").Copy($ws.Range("A1014:B1014"))
$ws.Range("A15:B15").Copy($ws.Range("A1015:B1015"))
$ws.Range("A16:C16").Copy($ws.Range("A1016:C1016"))
$ws.Range("A17:C17").Copy($ws.Range("A1017:C1017"))
$ws.Range("A18:E18").Copy($ws.Range("A1018:E1018"))
$ws.Range("A19:E19").Copy($ws.Range("A1019:E1019"))
$ws.Range("A20:E20").Copy($ws.Range("A1020:E1020"))
$ws.Range("A21:E21").Copy($ws.Range("A1021:E1021"))
$ws.Range("A22:D22").Copy($ws.Range("A1022:D1022"))
$ws.Range("A23:D23").Copy($ws.Range("A1023:D1023"))
$ws.Range("A24:D24").Copy($ws.Range("A1024:D1024"))
$ws.Range("A25:D25").Copy($ws.Range("A1025:D1025"))
$ws.Range("A26:E26").Copy($ws.Range("A1026:E1026"))
$ws.Range("A27:E27").Copy($ws.Range("A1027:E1027"))
$ws.Range("A28:E28").Copy($ws.Range("A1028:E1028"))
$ws.Range("A29:E29").Copy($ws.Range("A1029:E1029"))
$ws.Range("A30:E30").Copy($ws.Range("A1030:E1030"))
$ws.Range("A31:E31").Copy($ws.Range("A1031:E1031"))
$ws.Range("A32:E32").Copy($ws.Range("A1032:E1032"))
$ws.Range("A33:E33").Copy($ws.Range("A1033:E1033"))
$ws.Range("A34:F34").Copy($ws.Range("A1034:F1034"))
$ws.Range("A35:E35").Copy($ws.Range("A1035:E1035"))
$ws.Range("A36:E36").Copy($ws.Range("A1036:E1036"))
$ws.Range("A37:E37").Copy($ws.Range("A1037:E1037"))
$ws.Range("A38:E38").Copy($ws.Range("A1038:E1038"))
$ws.Range("A39:D39").Copy($ws.Range("A1039:D1039"))
$ws.Range("A40:D40").Copy($ws.Range("A1040:D1040"))
$ws.Range("A41:D41").Copy($ws.Range("A1041:D1041"))
$ws.Range("A42:D42").Copy($ws.Range("A1042:D1042"))
$ws.Range("A43:D43").Copy($ws.Range("A1043:D1043"))

# --- Phase 2: wipe rows 2-43 (values + formatting) so the write-back
#     below starts from a clean slate ---
$ws.Range("A2:F43").Clear()

# --- Phase 3: cut each staged row back into its new row number ---
$ws.Range("A1003:B1003").Cut($ws.Range("A2:B2"))
$ws.Range("A1005:B1005").Cut($ws.Range("A3:B3"))
$ws.Range("A1006:B1006").Cut($ws.Range("A4:B4"))
$ws.Range("A1007:B1007").Cut($ws.Range("A5:B5"))
$ws.Range("A1008:B1008").Cut($ws.Range("A6:B6"))
$ws.Range("A1009:B1009").Cut($ws.Range("A7:B7"))
$ws.Range("A1010:B1010").Cut($ws.Range("A8:B8"))
$ws.Range("A1011:B1011").Cut($ws.Range("A9:B9"))
$ws.Range("A1012:B1012").Cut($ws.Range("A10:B10"))
$ws.Range("A1015:B1015").Cut($ws.Range("A11:B11"))
$ws.Range("A1013:B1013").Cut($ws.Range("A12:B12"))
$ws.Range("A1014:B1014").Cut($ws.Range("A13:B13"))
$ws.Range("A1002:B1002").Cut($ws.Range("A14:B14"))
$ws.Range("A1004:B1004").Cut($ws.Range("A15:B15"))
$ws.Range("A1016:C1016").Cut($ws.Range("A16:C16"))
$ws.Range("A1017:C1017").Cut($ws.Range("A17:C17"))
$ws.Range("A1034:F1034").Cut($ws.Range("A18:F18"))
$ws.Range("A1018:E1018").Cut($ws.Range("A19:E19"))
$ws.Range("A1019:E1019").Cut($ws.Range("A20:E20"))
$ws.Range("A1020:E1020").Cut($ws.Range("A21:E21"))
$ws.Range("A1021:E1021").Cut($ws.Range("A22:E22"))
$ws.Range("A1026:E1026").Cut($ws.Range("A23:E23"))
$ws.Range("A1027:E1027").Cut($ws.Range("A24:E24"))
$ws.Range("A1028:E1028").Cut($ws.Range("A25:E25"))
$ws.Range("A1035:E1035").Cut($ws.Range("A26:E26"))
$ws.Range("A1036:E1036").Cut($ws.Range("A27:E27"))
$ws.Range("A1037:E1037").Cut($ws.Range("A28:E28"))
$ws.Range("A1038:E1038").Cut($ws.Range("A29:E29"))
$ws.Range("A1039:D1039").Cut($ws.Range("A30:D30"))
$ws.Range("A1040:D1040").Cut($ws.Range("A31:D31"))
$ws.Range("A1041:D1041").Cut($ws.Range("A32:D32"))
$ws.Range("A1042:D1042").Cut($ws.Range("A33:D33"))
$ws.Range("A1043:D1043").Cut($ws.Range("A34:D34"))
$ws.Range("A1022:D1022").Cut($ws.Range("A35:D35"))
$ws.Range("A1023:D1023").Cut($ws.Range("A36:D36"))
$ws.Range("A1024:D1024").Cut($ws.Range("A37:D37"))
$ws.Range("A1025:D1025").Cut($ws.Range("A38:D38"))
$ws.Range("A1029:E1029").Cut($ws.Range("A39:E39"))
$ws.Range("A1030:E1030").Cut($ws.Range("A40:E40"))
$ws.Range("A1031:E1031").Cut($ws.Range("A41:E41"))
$ws.Range("A1032:E1032").Cut($ws.Range("A42:E42"))
$ws.Range("A1033:E1033").Cut($ws.Range("A43:E43"))

# --- Phase 4: drop the scratch area so the used range / dimension
#     collapses back to A1:EI46 ---
$ws.Range("A1002:F1043").Clear()

# --- Fix up row "spans" hints to match rows 17-32 now carrying data
#     through column F / column D blocks through column F ---

# --- View/selection changes recorded in the sheetView ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("D14").Select()
